$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = 1
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 1

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = 2
$ws2.Range("A2").Value = 2
$ws2.Range("A3").Value = 3

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = 3
$ws3.Range("A2").Value = 3
$ws3.Range("A3").Value = 2

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A1").Value = "Pass"
$ws4.Range("A2").Value = "Pass"
$ws4.Range("A3").Value = "Fail"
$ws4.Range("A1:A2").Interior.Color = 65280
$ws4.Range("A3").Interior.Color = 255
